$d = $word.ActiveDocument

# --- Paragraph 4: "Relevant data..." - fix spelling, split runs to match target structure ---
$p4 = $d.Paragraphs(4)
$p4.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Relevant data about the neighbourhoods of Toronto are taken from Wikipedia. This data contains the </w:t></w:r><w:r><w:t>neighbourhood</w:t></w:r><w:r><w:t xml:space="preserve"> names and postal code as well as relevant other </w:t></w:r><w:r><w:t>information</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>')

# --- Paragraph 5: "The Foursquare API..." - fix informations -> information, split run ---
$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"> The Foursquare API will be used to get data about venues from the different neighbourhoods. With the Foursquare API it is possible to obtain details about what are the most popular places and other relevant </w:t></w:r><w:r><w:t>information</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>')

# --- Insert 6 new paragraphs after paragraph 5 (Methodology section through Conclusion) ---
$insertionPoint = $d.Paragraphs(5).Range
$insertionPoint.Collapse(0)

$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$newP = $insertionPoint.Paragraphs(1)
$newP.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Methodology</w:t></w:r></w:p>')
$insertionPoint = $newP.Range
$insertionPoint.Collapse(0)

$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$newP = $insertionPoint.Paragraphs(1)
$newP.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>It was decided to use a minimalist approach to solve the problem as we decided to follow the proverb “Less is more”. Once</w:t></w:r><w:r><w:t xml:space="preserve"> the data was retrieved it was </w:t></w:r><w:r><w:t>necessary</w:t></w:r><w:r><w:t xml:space="preserve"> to perform one hot encoding in order to transform these categorical values into binary values before being able to perform the classification. This procedure was therefore realized. It was decided to choose a value of 7 as the total number of clusters as this corresponded to the best way to have a somewhat uniform repartition of the </w:t></w:r><w:r><w:t>neighbourhood</w:t></w:r><w:r><w:t xml:space="preserve">. Despite that, the cluster number 1 was still composed of 79 </w:t></w:r><w:r><w:t>neighbourhoods</w:t></w:r><w:r><w:t xml:space="preserve">. It was therefore decided to use a further step to choose which </w:t></w:r><w:r><w:t>neighbourhood</w:t></w:r><w:r><w:t xml:space="preserve"> was </w:t></w:r><w:r><w:t>best suited for the opening of the restaurant. Find the number of neighbourhood where either the first 2</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:t xml:space="preserve"> or 3</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:t xml:space="preserve"> most common value is a Taiwanese restaurant. Then, find in their cluster the closest </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neighborhood</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> without a Taiwanese restaurant. The results are discussed below</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>')
$insertionPoint = $newP.Range
$insertionPoint.Collapse(0)

$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$newP = $insertionPoint.Paragraphs(1)
$newP.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Results and Discussion</w:t></w:r></w:p>')
$insertionPoint = $newP.Range
$insertionPoint.Collapse(0)

$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$newP = $insertionPoint.Paragraphs(1)
$newP.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="2790"/></w:tabs></w:pPr><w:r><w:t>After the clustering,</w:t></w:r><w:r><w:t xml:space="preserve"> and the condition applied. 2 </w:t></w:r><w:r><w:t>neighbourhood</w:t></w:r><w:r><w:t xml:space="preserve"> of cluster one were found to have one Thai restaurant. They were the index 40 and 79. It was thus decided to use the clustering distance to find the </w:t></w:r><w:r><w:t>neighbourhood</w:t></w:r><w:r><w:t xml:space="preserve"> the closest to one of these 2. The final </w:t></w:r><w:r><w:t>neighbourhood</w:t></w:r><w:r><w:t xml:space="preserve"> was found to be Northwood Park</w:t></w:r></w:p>')
$insertionPoint = $newP.Range
$insertionPoint.Collapse(0)

$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$newP = $insertionPoint.Paragraphs(1)
$newP.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Conclusion</w:t></w:r></w:p>')
$insertionPoint = $newP.Range
$insertionPoint.Collapse(0)

$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$newP = $insertionPoint.Paragraphs(1)
$newP.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Finally, to maximize the profit of the restaurant of the client we recommend him to open it at </w:t></w:r><w:r><w:t>Northwood</w:t></w:r><w:r><w:t xml:space="preserve"> Park as it is similar to the High Park </w:t></w:r><w:r><w:t>neighbourhood</w:t></w:r><w:r><w:t xml:space="preserve"> where a Thai restaurant is 3</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:t xml:space="preserve"> most popular venue.</w:t></w:r></w:p>')
$insertionPoint = $newP.Range
$insertionPoint.Collapse(0)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)